$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / comment corrections ---

# Heated bed comment: remove "might " before "get more printing area"
$ws.Range("E9").Value = "Optional, you can still use the printer without one. You get more printing area, but might not be able to print with abs."

# Threaded rod item name & spec update
$ws.Range("B12").Value = "5/16 - 18 Stainless Steel Threaded Rod (340mm)"
$ws.Range("E12").Value = "Z axis Threaded rod, if you want to use a Metric one, you can use a M8 Stainless Steel threaded rod. "

# Nuts item name & comment update
$ws.Range("B13").Value = "5/16 - 18 Stainless Steel Nuts"
$ws.Range("E13").Value = "You only need one, more just in case. For Metric, you can use an M8 nut."

# M3 Washers row: add a new comment explaining their purpose
$ws.Range("E32").Value = "Washers to get good spacing between parts. "

# --- Column width tweak (Comment/Source column D got a bit wider) ---
$ws.Range("D:D").ColumnWidth = 36.309523809524

# --- Update view: scroll/selection moved to the row that was edited ---
$ws.Range("E9").Select()
